$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the existing activity descriptions with more detail
$ws.Range("A2").Value = "Desingning Forms(Login, register, main menu)"
$ws.Range("A3").Value = "Creating Fomrs (Login, register, main menu)"

# Add a new activity row (row 4)
# Dates: copy the date formatting used by the row above so the new cells
# keep the same numeric date style (B column), and a right-aligned date
# style for column C.
$ws.Range("B3").Copy($ws.Range("B4"))
$ws.Range("B4").Value = 43738

$ws.Range("C3").Copy($ws.Range("C4"))
$ws.Range("C4").Value = 43738
$ws.Range("C4").HorizontalAlignment = -4152

$ws.Range("D4").Value = "5pm"
$ws.Range("E4").Value = "11pm"
$ws.Range("A4").Value = "Creating and Designing of Forms (Check-In, Check-Out, Reserve a Room)"

# Widen column A so the longer activity text fits
$ws.Columns.Item(1).ColumnWidth = 67.7

$null = $ws.Range("E5").Select()
